$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.728999999999999
$ws.Range("C3").Value = -11.425
$ws.Range("D3").Value = -7.281999999999999
$ws.Range("C4").Value = -12.943
$ws.Range("D9").Value = -7.162000000000001
$ws.Range("A11").Value = -21.81
$ws.Range("A12").Value = -21.669
$ws.Range("C14").Value = -12.543
$ws.Range("A15").Value = -22.015
$ws.Range("D15").Value = -8.016000000000002
$ws.Range("D19").Value = -8.135999999999999
$ws.Range("D20").Value = -7.825
$ws.Range("D25").Value = -7.906000000000001
$ws.Range("C26").Value = -12.868
$ws.Range("A27").Value = -21.749
$ws.Range("D27").Value = -8.119
$ws.Range("A28").Value = -21.875
$ws.Range("D28").Value = -7.85
$ws.Range("D30").Value = -7.202
$ws.Range("A31").Value = -21.96
$ws.Range("C31").Value = -13.342
$ws.Range("A32").Value = -21.717
$ws.Range("D32").Value = -7.739999999999999
$ws.Range("C35").Value = -12.288
$ws.Range("A36").Value = -20.257
$ws.Range("C37").Value = -13.38
$ws.Range("A38").Value = -19.696
$ws.Range("C39").Value = -12.745
$ws.Range("C40").Value = -12.989
$ws.Range("D44").Value = -7.904999999999999
$ws.Range("C45").Value = -12.569
$ws.Range("A46").Value = -21.883
$ws.Range("D47").Value = -7.367
$ws.Range("C52").Value = -11.363
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.21
$ws.Range("A56").Value = -21.997
$ws.Range("C57").Value = -13.337
$ws.Range("D58").Value = -8.218
$ws.Range("D62").Value = -7.994
$ws.Range("A67").Value = -21.517
$ws.Range("A69").Value = -21.637
$ws.Range("A72").Value = -21.481
$ws.Range("A73").Value = -19.973
$ws.Range("D77").Value = -7.949000000000001
$ws.Range("D78").Value = -7.972999999999999
$ws.Range("C81").Value = -13.335
$ws.Range("A83").Value = -21.628
$ws.Range("C83").Value = -12.893
$ws.Range("D84").Value = -8.019
$ws.Range("A86").Value = -22.257
$ws.Range("D89").Value = -6.946
$ws.Range("A91").Value = -21.584
$ws.Range("D91").Value = -7.053
$ws.Range("D92").Value = -6.746
$ws.Range("A93").Value = -21.547
$ws.Range("D96").Value = -7.297
$ws.Range("A99").Value = -19.854
$ws.Range("C100").Value = -12.638
$ws.Range("C102").Value = -13.281
$ws.Range("D102").Value = -8.026999999999999
